$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.471.93"
$ws.Range("E2").Value = "  -0.98%  "
$ws.Range("D3").Value = "3.897.57"
$ws.Range("E3").Value = "  +2.72%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'602.11"
$ws.Range("E5").Value = "  +0.08%  "
$ws.Range("D6").Value = "'166.26"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").Value = "3.894.77"
$ws.Range("E7").Value = "  +2.69%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  -1.27%  "
$ws.Range("D11").Value = "'6.41"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D13").Value = "'0.0000255"
$ws.Range("E13").Value = "  +3.81%  "
$ws.Range("D14").Value = "'37.31"
$ws.Range("E14").Value = "  +0.08%  "
$ws.Range("D15").Value = "4.550.41"
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("D16").Value = "3.898.71"
$ws.Range("E16").Value = "  +2.87%  "
$ws.Range("D17").Value = "68.552.00"
$ws.Range("E17").Value = "  -0.97%  "
$ws.Range("D19").Value = "'17.16"
$ws.Range("E19").Value = "  -0.66%  "
$ws.Range("E20").Value = "  -2.24%  "
$ws.Range("D21").Value = "'11.03"
$ws.Range("E21").Value = "  -2.96%  "
$ws.Range("D22").Value = "'487.10"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("E24").Value = "  +10.51%  "
$ws.Range("E25").Value = "  -0.14%  "
$ws.Range("D26").Value = "'2.24"
$ws.Range("E26").Value = "  -0.82%  "
$ws.Range("D27").Value = "'12.06"
$ws.Range("E27").Value = "  -1.33%  "
$ws.Range("D28").Value = "'10.12"
$ws.Range("E28").Value = "  +0.77%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -0.82%  "
$ws.Range("D31").Value = "4.048.90"
$ws.Range("E31").Value = "  +2.62%  "
$ws.Range("E32").Value = "  -0.59%  "
$ws.Range("E33").Value = "  -3.87%  "
$ws.Range("E34").Value = "  +0.20%  "
$ws.Range("D35").Value = "3.850.84"
$ws.Range("E35").Value = "  +2.91%  "
$ws.Range("E36").Value = "  +0.02%  "
$ws.Range("E37").Value = "  +2.15%  "
$ws.Range("D38").Value = "'5.94"
$ws.Range("E38").Value = "  +0.47%  "
$ws.Range("E39").Value = "  -1.96%  "
$ws.Range("E40").Value = "  +5.40%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -2.07%  "
$ws.Range("D43").Value = "'429.49"
$ws.Range("E43").Value = "  +1.93%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.99"
$ws.Range("E44").Value = "  -0.11%  "
$ws.Range("B45").Value = "OKB"
$ws.Range("C45").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D45").Value = "'48.31"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("E47").Value = "  -0.01%  "
$ws.Range("D48").Value = "'142.20"
$ws.Range("E48").Value = "  +0.74%  "
$ws.Range("D49").Value = "'26.07"
$ws.Range("E49").Value = "  +8.58%  "
$ws.Range("D50").Value = "2.808.50"
$ws.Range("E50").Value = "  -0.89%  "
$ws.Range("E51").Value = "  +0.67%  "
